$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark from its original location
#    (it previously sat in an empty paragraph right after "$SZ").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the paragraph that lists the special circumstances
#    "(z.B. Nachholfristen, Gastschüler/innen, Probezeit, freiwilliger
#    Rücktritt)" and append a new run ", besonders begabte
#    Schüler/innen" right before the closing parenthesis, then place
#    the (now recreated) "_GoBack" bookmark immediately after that
#    new run.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("(z.B. Nachholfristen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $rng.Paragraphs(1).Range

    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00FF0B22" w:rsidRPr="005A3C83" w:rsidRDefault="00FF0B22" w:rsidP="005A3C83"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:spacing w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="de-DE"/></w:rPr></w:pPr><w:r w:rsidRPr="005A3C83"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:spacing w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>(z.B. Nachholfristen, Gastsch&#252;ler/innen, Probezeit</w:t></w:r><w:r w:rsidR="005A3C83"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:spacing w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>, freiwilliger R&#252;cktritt</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:spacing w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>, besonders begabte Sch&#252;ler/innen</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="005A3C83"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:spacing w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="de-DE"/></w:rPr><w:t>)</w:t></w:r><w:r w:rsidRPr="005A3C83"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:spacing w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="de-DE"/></w:rPr><w:br/></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

    $para.InsertXML($xml)
}
